$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds plain-text figures (e.g. "1.002", "22.467.60")
# that look numeric. Mark the range as Text first so Excel's normal
# type-inference does not silently convert them to numbers/dates when
# the new values are written below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.467.60"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "1.573.94"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "291.54"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").Value = "0.3749"
$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").Value = "49.97"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "0.3410"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").Value = "1.151"
$ws.Range("E10").Value = "  -1.23%  "

$ws.Range("D11").Value = "0.07586"
$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").Value = "21.46"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").Value = "6.004"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").Value = "6.971"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "1.572.05"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "0.00001126"
$ws.Range("E17").Value = "  -0.86%  "

$ws.Range("D18").Value = "91.35"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").Value = "0.06741"
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").Value = "6.289"
$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").Value = "16.46"
$ws.Range("E22").Value = "  -1.75%  "

$ws.Range("D23").Value = "12.20"
$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("D24").Value = "22.466.70"
$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("D25").Value = "2.328"
$ws.Range("E25").Value = "  -3.70%  "

$ws.Range("D26").Value = "2.608"
$ws.Range("E26").Value = "  -5.25%  "

$ws.Range("D27").Value = "20.17"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").Value = "148.69"
$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("D29").Value = "5.009"
$ws.Range("E29").Value = "  -1.05%  "

$ws.Range("D30").Value = "126.27"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "1.747.94"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").Value = "1.043"
$ws.Range("E32").Value = "  +2.65%  "

$ws.Range("D33").Value = "6.165"
$ws.Range("E33").Value = "  -0.94%  "

$ws.Range("E34").Value = "  -2.13%  "

$ws.Range("D35").Value = "9.921"
$ws.Range("E35").Value = "  -1.79%  "

$ws.Range("D36").Value = "0.08474"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("D37").Value = "1.391"
$ws.Range("E37").Value = "  +4.68%  "

$ws.Range("D38").Value = "0.02474"
$ws.Range("E38").Value = "  -3.42%  "

$ws.Range("D39").Value = "0.2300"
$ws.Range("E39").Value = "  -0.81%  "

$ws.Range("D40").Value = "0.06554"
$ws.Range("E40").Value = "  +0.38%  "

$ws.Range("D41").Value = "5.510"
$ws.Range("E41").Value = "  +0.56%  "

$ws.Range("D42").Value = "11.44"
$ws.Range("E42").Value = "  -1.67%  "

$ws.Range("D43").Value = "0.6303"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").Value = "14.08"
$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").Value = "3.816"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("D47").Value = "0.5898"
$ws.Range("E47").Value = "  -2.39%  "

$ws.Range("D48").Value = "2.101"
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("D49").Value = "130.56"
$ws.Range("E49").Value = "  +3.85%  "

$ws.Range("E50").Value = "  -6.12%  "

$ws.Range("D51").Value = "0.07344"
$ws.Range("E51").Value = "  -0.05%  "
